$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Price / percentage columns are stored as plain text in this sheet
    # (e.g. "294.80", "1.96%"). A bare Range.Value assignment lets Excel's
    # automatic type detection turn these into numbers/percentages, which
    # loses the original text formatting (trailing zeros, "%" literal,
    # leading zeros, etc). Force the cell to Text first, write the literal
    # string, then drop back to the default (General) number format so the
    # cell's style matches the un-styled original cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# --- Rows 6-17: "Coin" (B) / "Link" (C) cycle up by one row (row 6 picks up
#     what used to be in row 7, ... row 17 wraps around to what used to be
#     in row 6), and the Price (D) / Volume(1h) (E) figures refresh to the
#     newly pulled values. ---

$coinRows = @(
    @{ Row = 6;  B = "FTXToken";                         C = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt";                                   D = "1.571";   E = "1.01%" }
    @{ Row = 7;  B = "MXToken";                           C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx";                                   D = "0.9314";  E = "2.31%" }
    @{ Row = 8;  B = "BTSEToken";                          C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";                               D = "2.394";   E = "-0.12%" }
    @{ Row = 9;  B = "LiechtensteinCryptoassetsExchange"; C = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx";            D = "0.1193";  E = "-0.02%" }
    @{ Row = 10; B = "WazirX";                             C = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx";                                        D = "0.1801";  E = "3.44%" }
    @{ Row = 11; B = "MandalaExchangeToken";               C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx";                     D = "0.08743"; E = "0.68%" }
    @{ Row = 12; B = "BitrueCoin";                         C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr";                                   D = "0.04274"; E = "2.29%" }
    @{ Row = 13; B = "BitMartToken";                       C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";                             D = "0.1053";  E = "0.19%" }
    @{ Row = 14; B = "BitForexToken";                      C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";                             D = "0.001266"; E = "-0.73%" }
    @{ Row = 15; B = "TigerCash";                          C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                                     D = "0.005947"; E = "0.55%" }
    @{ Row = 16; B = "LEO";                                C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                                       D = "3.341";   E = "-1.36%" }
    @{ Row = 17; B = "GateToken";                          C = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt";                                  D = "4.323";   E = "0.86%" }
)

foreach ($entry in $coinRows) {
    $ws.Range("B$($entry.Row)").Value = $entry.B
    $ws.Range("C$($entry.Row)").Value = $entry.C
    Set-TextValue "D$($entry.Row)" $entry.D
    Set-TextValue "E$($entry.Row)" $entry.E
}

# --- Remaining rows: only Price (D) and/or Volume(1h) (E) refreshed ---

$valueRows = @(
    @{ Row = 2;  D = "296.71";        E = "1.93%" }
    @{ Row = 3;  D = "40.68";         E = "1.88%" }
    @{ Row = 4;  D = "5.015";         E = "-0.02%" }
    @{ Row = 5;  D = "0.07440";       E = "0.95%" }
    @{ Row = 18; D = "0.3297";        E = "-0.13%" }
    @{ Row = 19; D = "7.867";         E = "4.41%" }
    @{ Row = 20; D = "0.1379";        E = "1.95%" }
    @{ Row = 21; D = "0.3294";        E = "14.27%" }
    @{ Row = 22; D = "0.03947";       E = "2.91%" }
    @{ Row = 23; E = "-0.40%" }
    @{ Row = 24; D = "0.003822";      E = "-1.83%" }
    @{ Row = 25; D = "0.0001222";     E = "-4.61%" }
    @{ Row = 26; D = "0.0003714";     E = "-0.28%" }
    @{ Row = 38; D = "0.02360";       E = "1.10%" }
    @{ Row = 39; D = "0.05127";       E = "1.98%" }
    @{ Row = 40; D = "0.005900";      E = "15.47%" }
    @{ Row = 41; D = "0.007727";      E = "0.40%" }
    @{ Row = 42; D = "0.1314";        E = "3.44%" }
    @{ Row = 43; D = "0.007366";      E = "-0.13%" }
    @{ Row = 44; D = "0.007001";      E = "0.51%" }
    @{ Row = 45; D = "0.2944";        E = "-6.68%" }
    @{ Row = 46; D = "0.00006205";    E = "-4.67%" }
    @{ Row = 47; D = "0.00000000745"; E = "-0.70%" }
    @{ Row = 48; D = "0.04630";       E = "-81.61%" }
    @{ Row = 49; D = "0.004191";      E = "-0.27%" }
    @{ Row = 50; D = "0.00002086";    E = "-0.70%" }
    @{ Row = 51; D = "0.0001987";     E = "-0.70%" }
)

foreach ($entry in $valueRows) {
    if ($entry.ContainsKey("D")) {
        Set-TextValue "D$($entry.Row)" $entry.D
    }
    if ($entry.ContainsKey("E")) {
        Set-TextValue "E$($entry.Row)" $entry.E
    }
}
